$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the missing Form name for CO and update its Margin value
$ws.Range("B6").Value = "CO DR-0100 XML"
$ws.Range("C6").Value = 1

# Remove the stray leftover row 21 ("Questions" in column A only)
$ws.Rows.Item(21).Delete()

$ws.Range("H10").Select()
